$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$wsMetadata = $wb.Worksheets.Item("Metadata")

# Date: 2026-01-09T15:21:06+00:00 -> 2026-01-16T13:49:34+00:00
$wsMetadata.Range("B8").Value = "2026-01-16T13:49:34+00:00"

# Description: "Entrée Resultats d'examens de biologie medicale" -> "Resultats d'examens de biologie medicale"
$wsMetadata.Range("B12").Value = "Resultats d'examens de biologie medicale"

# --- Elements sheet updates ---
$wsElements = $wb.Worksheets.Item("Elements")

# Short/Definition for the root element (row 2) mirrors the Description change
$wsElements.Range("M2").Value = "Resultats d'examens de biologie medicale"

# laboratoireExecutant (row 7): trim the trailing explanatory sentence
$wsElements.Range("L7").Value = "Laboratoire sous-traitant."
$wsElements.Range("M7").Value = "Laboratoire sous-traitant."

# auteur (row 8): trim the trailing explanatory sentence
$wsElements.Range("L8").Value = "Participation d'un auteur au document."
$wsElements.Range("M8").Value = "Participation d'un auteur au document."
